# Scheduled-runner refresh of market/profit figures (columns H:N) across
# all item-category sheets. Values below are the refreshed snapshot
# (currentAveragePrice*, LevePrice*, LeveProfit*) recomputed upstream;
# we just write the new numbers into the same cells.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 366.66666
$ws.Range("I2").Value = 375
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 375
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -262
$ws.Range("N2").Value = -576

$ws.Range("H74").Value = 4718.6875
$ws.Range("I74").Value = 4283.8335
$ws.Range("J74").Value = 4979.6
$ws.Range("K74").Value = 4283.8335
$ws.Range("L74").Value = 4979.6
$ws.Range("M74").Value = -3347.8335
$ws.Range("N74").Value = -6851.6

$ws.Range("H76").Value = 6817.1665
$ws.Range("J76").Value = 7900
$ws.Range("L76").Value = 7900
$ws.Range("N76").Value = -8530

$ws.Range("H77").Value = 4718.6875
$ws.Range("I77").Value = 4283.8335
$ws.Range("J77").Value = 4979.6
$ws.Range("K77").Value = 21419.1675
$ws.Range("L77").Value = 24898
$ws.Range("M77").Value = -16739.1675
$ws.Range("N77").Value = -34258

$ws.Range("H79").Value = 6817.1665
$ws.Range("J79").Value = 7900
$ws.Range("L79").Value = 7900
$ws.Range("N79").Value = -10084

$ws.Range("H81").Value = 35998.332
$ws.Range("J81").Value = 35998.332
$ws.Range("L81").Value = 35998.332
$ws.Range("N81").Value = -37994.332

$ws.Range("H84").Value = 35998.332
$ws.Range("J84").Value = 35998.332
$ws.Range("L84").Value = 107994.996
$ws.Range("N84").Value = -117978.996

$ws.Range("H121").Value = 1748.6666
$ws.Range("I121").Value = 876
$ws.Range("J121").Value = 2185
$ws.Range("K121").Value = 2628
$ws.Range("L121").Value = 6555
$ws.Range("M121").Value = -881
$ws.Range("N121").Value = -10049

$ws.Range("H129").Value = 1045.878
$ws.Range("I129").Value = 365.57144
$ws.Range("J129").Value = 1398.6296
$ws.Range("K129").Value = 1096.71432
$ws.Range("L129").Value = 4195.8888
$ws.Range("M129").Value = 3903.28568
$ws.Range("N129").Value = -14195.8888

$ws.Range("H131").Value = 6411.3
$ws.Range("I131").Value = 922.9
$ws.Range("K131").Value = 2768.7
$ws.Range("M131").Value = 2271.3

$ws.Range("H138").Value = 1421.41
$ws.Range("I138").Value = 786.6896400000001
$ws.Range("J138").Value = 2297.9285
$ws.Range("K138").Value = 2360.06892
$ws.Range("L138").Value = 6893.7855
$ws.Range("M138").Value = 2779.93108
$ws.Range("N138").Value = -17173.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14762.357
$ws.Range("I32").Value = 15465.583
$ws.Range("K32").Value = 15465.583
$ws.Range("M32").Value = -15178.583

$ws.Range("H132").Value = 6107.3438
$ws.Range("I132").Value = 8636.647000000001
$ws.Range("K132").Value = 25909.941
$ws.Range("M132").Value = -23379.941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 29639.334
$ws.Range("J62").Value = 30327.857
$ws.Range("L62").Value = 30327.857
$ws.Range("N62").Value = -31699.857

$ws.Range("H65").Value = 29639.334
$ws.Range("J65").Value = 30327.857
$ws.Range("L65").Value = 90983.571
$ws.Range("N65").Value = -97847.571

$ws.Range("H80").Value = 2314970
$ws.Range("I80").Value = 12345715
$ws.Range("J80").Value = 182.53847
$ws.Range("K80").Value = 12345715
$ws.Range("L80").Value = 182.53847
$ws.Range("M80").Value = -12344717
$ws.Range("N80").Value = -2178.53847

$ws.Range("H83").Value = 2314970
$ws.Range("I83").Value = 12345715
$ws.Range("J83").Value = 182.53847
$ws.Range("K83").Value = 61728575
$ws.Range("L83").Value = 912.6923499999999
$ws.Range("M83").Value = -61723583
$ws.Range("N83").Value = -10896.69235

$ws.Range("H86").Value = 503253.25
$ws.Range("I86").Value = 10006
$ws.Range("J86").Value = 667669
$ws.Range("K86").Value = 10006
$ws.Range("L86").Value = 667669
$ws.Range("M86").Value = -8883
$ws.Range("N86").Value = -669915

$ws.Range("H89").Value = 503253.25
$ws.Range("I89").Value = 10006
$ws.Range("J89").Value = 667669
$ws.Range("K89").Value = 50030
$ws.Range("L89").Value = 3338345
$ws.Range("M89").Value = -44414
$ws.Range("N89").Value = -3349577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1776
$ws.Range("I122").Value = 1726.48
$ws.Range("J122").Value = 3014
$ws.Range("K122").Value = 5179.440000000001
$ws.Range("L122").Value = 9042
$ws.Range("M122").Value = -2729.440000000001
$ws.Range("N122").Value = -13942

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1780.0588
$ws.Range("I5").Value = 1935.0769
$ws.Range("J5").Value = 1276.25
$ws.Range("K5").Value = 5805.2307
$ws.Range("L5").Value = 3828.75
$ws.Range("M5").Value = -5693.2307
$ws.Range("N5").Value = -4052.75

$ws.Range("H21").Value = 4631562.5
$ws.Range("I21").Value = 799.5
$ws.Range("J21").Value = 6946944
$ws.Range("K21").Value = 2398.5
$ws.Range("L21").Value = 20840832
$ws.Range("M21").Value = -2225.5
$ws.Range("N21").Value = -20841178

$ws.Range("H131").Value = 1629.2452
$ws.Range("J131").Value = 1784.8937
$ws.Range("L131").Value = 5354.6811
$ws.Range("N131").Value = -15434.6811

$ws.Range("H132").Value = 1942.7858
$ws.Range("J132").Value = 2704.3
$ws.Range("L132").Value = 24338.7
$ws.Range("N132").Value = -29398.7

$ws.Range("H135").Value = 1780.0588
$ws.Range("I135").Value = 1935.0769
$ws.Range("J135").Value = 1276.25
$ws.Range("K135").Value = 17415.6921
$ws.Range("L135").Value = 11486.25
$ws.Range("M135").Value = -14880.6921
$ws.Range("N135").Value = -16556.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7797.5127
$ws.Range("I70").Value = 7956.9287
$ws.Range("K70").Value = 7956.9287
$ws.Range("M70").Value = -7686.9287

$ws.Range("H73").Value = 7797.5127
$ws.Range("I73").Value = 7956.9287
$ws.Range("K73").Value = 7956.9287
$ws.Range("M73").Value = -7020.9287

$ws.Range("J80").Value = 3250
$ws.Range("L80").Value = 3250
$ws.Range("N80").Value = -5246

$ws.Range("J83").Value = 3250
$ws.Range("L83").Value = 16250
$ws.Range("N83").Value = -26234

$ws.Range("H109").Value = 9270.875
$ws.Range("J109").Value = 9270.875
$ws.Range("L109").Value = 9270.875
$ws.Range("N109").Value = -11350.875

$ws.Range("H132").Value = 2673.0889
$ws.Range("I132").Value = 2283.6875
$ws.Range("J132").Value = 3631.6155
$ws.Range("K132").Value = 6851.0625
$ws.Range("L132").Value = 10894.8465
$ws.Range("M132").Value = -4321.0625
$ws.Range("N132").Value = -15954.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6000
$ws.Range("I7").Value = 6666.6665
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 6666.6665
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -6554.6665
$ws.Range("N7").Value = -4224

$ws.Range("H40").Value = 4540
$ws.Range("I40").Value = 5297.143
$ws.Range("K40").Value = 5297.143
$ws.Range("M40").Value = -5161.143

$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 1625
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 1625
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -2001

$ws.Range("H82").Value = 1820
$ws.Range("I82").Value = 1820
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1820
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -1459

$ws.Range("H85").Value = 1820
$ws.Range("I85").Value = 1820
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1820
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -572

$ws.Range("H122").Value = 34621276
$ws.Range("I122").Value = 83338664
$ws.Range("K122").Value = 250015992
$ws.Range("M122").Value = -250013542

$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 6666.6665
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 19999.9995
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -17529.9995
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 2554.0894
$ws.Range("I132").Value = 1615.0294
$ws.Range("J132").Value = 4005.3635
$ws.Range("K132").Value = 4845.0882
$ws.Range("L132").Value = 12016.0905
$ws.Range("M132").Value = -2315.0882
$ws.Range("N132").Value = -17076.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4644.4443

$ws.Range("H65").Value = 4644.4443

$ws.Range("H122").Value = 17363962
$ws.Range("I122").Value = 27780232
$ws.Range("J122").Value = 6947690.5
$ws.Range("K122").Value = 83340696
$ws.Range("L122").Value = 20843071.5
$ws.Range("M122").Value = -83338246
$ws.Range("N122").Value = -20847971.5

$ws.Range("H123").Value = 24009.064
$ws.Range("J123").Value = 24009.064
$ws.Range("L123").Value = 24009.064
$ws.Range("N123").Value = -33809.064
